# daily auto push: 2026-02-27 19:00 UTC
# Insert a new data row at row 899 (pushes existing rows 899-940 down to 900-941)
# and populate it with the new day's first reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(899).Insert()

# Column A holds dates stored as plain text (e.g. "2026/02/28"), not real
# date serials. Assigning that literal string makes Excel auto-detect it as
# a date, so force it to stay text the same way a user typing a leading
# apostrophe would, then drop the resulting "Text" number-format style so
# the cell matches its plain, unstyled neighbours.
$ws.Range("A899").Value = "'2026/02/28"
$ws.Range("A899").Style = "Normal"

$ws.Range("B899").Value = "土"
$ws.Range("C899").Value = 1
$ws.Range("D899").Value = 201
